$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.139.17'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.868.17'
$ws.Range("E3").Value = '  +3.82%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.59'
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("E7").Value = '  -1.59%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3904'
$ws.Range("E8").Value = '  +1.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09693'
$ws.Range("E9").Value = '  +25.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.137'
$ws.Range("E10").Value = '  +3.72%  '
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.456'
$ws.Range("E12").Value = '  +1.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.89'
$ws.Range("E13").Value = '  +2.95%  '
$ws.Range("D14").Value = '1.863.00'
$ws.Range("E14").Value = '  +3.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.9999'
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.383'
$ws.Range("E16").Value = '  +1.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001123'
$ws.Range("E17").Value = '  +5.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.04'
$ws.Range("E18").Value = '  +1.02%  '
$ws.Range("E19").Value = '  +0.59%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.42'
$ws.Range("E20").Value = '  +1.15%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.124'
$ws.Range("E22").Value = '  +2.81%  '
$ws.Range("D23").Value = '28.193.75'
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.32'
$ws.Range("E24").Value = '  +2.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.276'
$ws.Range("E25").Value = '  +1.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.559'
$ws.Range("E26").Value = '  +6.21%  '
$ws.Range("D27").Value = '2.084.34'
$ws.Range("E27").Value = '  +3.61%  '
$ws.Range("E28").Value = '  +4.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '157.21'
$ws.Range("E29").Value = '  -1.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.33'
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  -2.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.061'
$ws.Range("E32").Value = '  +1.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.630'
$ws.Range("E33").Value = '  +1.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.622'
$ws.Range("E34").Value = '  -0.81%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06753'
$ws.Range("E35").Value = '  -2.97%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.544'
$ws.Range("E36").Value = '  +5.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02394'
$ws.Range("E37").Value = '  +2.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2180'
$ws.Range("E38").Value = '  +0.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.48'
$ws.Range("E39").Value = '  +0.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.993'
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6284'
$ws.Range("E41").Value = '  +2.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.170'
$ws.Range("E42").Value = '  +1.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.47'
$ws.Range("E44").Value = '  +1.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6012'
$ws.Range("E45").Value = '  +2.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.664'
$ws.Range("E46").Value = '  -1.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.260'
$ws.Range("E47").Value = '  -2.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.21'
$ws.Range("E48").Value = '  -0.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.974'
$ws.Range("E49").Value = '  +2.92%  '
$ws.Range("E50").Value = '  +0.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06831'
$ws.Range("E51").Value = '  +1.46%  '
